$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in RANK formulas for K2:K11 based on Total Marks (column I)
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("K$r").Formula = "=RANK(I$r,`$I`$2:`$I`$11)"
}

# Update the active selection to match the target state
$ws.Range("K2:K11").Select()
